# Adds two new columns, I ("I0") and J ("IF"), to Sheet1.
# Header cells get the same bold/centered/bordered style as the
# existing header row (copied from H1), and each data row (2-68)
# gets its corresponding numeric value for I and J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, borders, center/top alignment)
# from the existing H1 header cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows ---
$data = @(
    @(2, 8, 8),
    @(3, 8, 8),
    @(4, 9, 9),
    @(5, 9, 9),
    @(6, 9, 9),
    @(7, 9, 9),
    @(8, 9, 9),
    @(9, 9, 9),
    @(10, 11, 11),
    @(11, 10, 10),
    @(12, 9, 9),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 8, 8),
    @(16, 9, 9),
    @(17, 7, 7),
    @(18, 9, 9),
    @(19, 9, 9),
    @(20, 8, 8),
    @(21, 9, 9),
    @(22, 9, 9),
    @(23, 10, 10),
    @(24, 9, 9),
    @(25, 9, 9),
    @(26, 9, 9),
    @(27, 8, 8),
    @(28, 9, 9),
    @(29, 9, 9),
    @(30, 8, 8),
    @(31, 9, 9),
    @(32, 9, 9),
    @(33, 9, 9),
    @(34, 9, 9),
    @(35, 9, 9),
    @(36, 9, 9),
    @(37, 9, 9),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 9, 9),
    @(41, 9, 9),
    @(42, 9, 9),
    @(43, 9, 9),
    @(44, 9, 10),
    @(45, 8, 8),
    @(46, 9, 9),
    @(47, 9, 9),
    @(48, 9, 9),
    @(49, 9, 9),
    @(50, 9, 9),
    @(51, 9, 9),
    @(52, 9, 9),
    @(53, 9, 9),
    @(54, 9, 9),
    @(55, 7, 8),
    @(56, 9, 9),
    @(57, 9, 9),
    @(58, 9, 9),
    @(59, 9, 9),
    @(60, 8, 9),
    @(61, 9, 9),
    @(62, 9, 9),
    @(63, 9, 9),
    @(64, 6, 7),
    @(65, 4, 4),
    @(66, 3, 3),
    @(67, 4, 4),
    @(68, 3, 3)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($row, 9).Value = $iVal   # column I
    $ws.Cells.Item($row, 10).Value = $jVal  # column J
}
